$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price (D) and Volume(1h) (E) columns are stored as plain text in the
# workbook, even when a price looks like a number (e.g. "0.997"). A literal
# leading apostrophe (the PowerShell literal '''...' yields a string that
# starts with a single quote) forces Excel to keep such values as text
# instead of silently converting them to numeric cells.

$ws.Range("D2").Value = '32.940.17'
$ws.Range("E2").Value = '  +10.13%  '
$ws.Range("D3").Value = '1.755.64'
$ws.Range("E3").Value = '  +5.90%  '
$ws.Range("D4").Value = '''0.997'
$ws.Range("E4").Value = '  -0.21%  '
$ws.Range("D5").Value = '''226.84'
$ws.Range("E5").Value = '  +4.50%  '
$ws.Range("E6").Value = '  +4.45%  '
$ws.Range("D7").Value = '''0.995'
$ws.Range("E7").Value = '  -0.36%  '
$ws.Range("E8").Value = '  +10.78%  '
$ws.Range("D9").Value = '''44.99'
$ws.Range("E9").Value = '  +2.71%  '
$ws.Range("D10").Value = '''0.279'
$ws.Range("E10").Value = '  +5.88%  '
$ws.Range("D11").Value = '''0.0664'
$ws.Range("E11").Value = '  +8.18%  '
$ws.Range("D12").Value = '''0.0918'
$ws.Range("E12").Value = '  +1.77%  '
$ws.Range("D13").Value = '2.008.15'
$ws.Range("E13").Value = '  +5.92%  '
$ws.Range("D14").Value = '1.762.27'
$ws.Range("E14").Value = '  +6.18%  '
$ws.Range("B15").Value = 'Polygon'
$ws.Range("C15").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D15").Value = '''0.631'
$ws.Range("E15").Value = '  +4.97%  '
$ws.Range("B16").Value = 'Chainlink'
$ws.Range("C16").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''10.51'
$ws.Range("E16").Value = '  +5.15%  '
$ws.Range("D17").Value = '''4.28'
$ws.Range("E17").Value = '  +8.68%  '
$ws.Range("D18").Value = '32.927.92'
$ws.Range("E18").Value = '  +9.98%  '
$ws.Range("D19").Value = '''68.75'
$ws.Range("E19").Value = '  +6.03%  '
$ws.Range("D20").Value = '''259.51'
$ws.Range("E20").Value = '  +7.25%  '
$ws.Range("D21").Value = '0.0₃0740'
$ws.Range("E21").Value = '  +4.39%  '
$ws.Range("D22").Value = '''0.999'
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = '''10.49'
$ws.Range("E23").Value = '  +4.77%  '
$ws.Range("E24").Value = '  +4.29%  '
$ws.Range("E25").Value = '  -0.85%  '
$ws.Range("D26").Value = '''159.98'
$ws.Range("E26").Value = '  +1.19%  '
$ws.Range("E27").Value = '  +5.11%  '
$ws.Range("E28").Value = '  +4.30%  '
$ws.Range("D29").Value = '''6.97'
$ws.Range("E29").Value = '  +3.75%  '
$ws.Range("D30").Value = '''0.995'
$ws.Range("E30").Value = '  -0.29%  '
$ws.Range("D31").Value = '''3.90'
$ws.Range("E31").Value = '  +14.76%  '
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("E33").Value = '  +5.54%  '
$ws.Range("D34").Value = '''3.47'
$ws.Range("E34").Value = '  +7.89%  '
$ws.Range("D35").Value = '1.553.84'
$ws.Range("E35").Value = '  +7.69%  '
$ws.Range("E36").Value = '  +4.94%  '
$ws.Range("B37").Value = 'ImmutableX'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D37").Value = '''0.631'
$ws.Range("E37").Value = '  +10.03%  '
$ws.Range("B38").Value = 'TrustWalletToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D38").Value = '''1.04'
$ws.Range("E38").Value = '  +1.53%  '
$ws.Range("B39").Value = 'Aave'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D39").Value = '''84.92'
$ws.Range("E39").Value = '  +7.73%  '
$ws.Range("D40").Value = '''0.0186'
$ws.Range("E40").Value = '  +6.09%  '
$ws.Range("E41").Value = '  +1.61%  '
$ws.Range("E42").Value = '  +0.11%  '
$ws.Range("D43").Value = '''2.09'
$ws.Range("E43").Value = '  +7.59%  '
$ws.Range("E44").Value = '  +3.40%  '
$ws.Range("D45").Value = '''0.0514'
$ws.Range("E45").Value = '  +2.56%  '
$ws.Range("D46").Value = '''55.04'
$ws.Range("E46").Value = '  +8.41%  '
$ws.Range("E47").Value = '  +4.47%  '
$ws.Range("E48").Value = '  +5.64%  '
$ws.Range("E49").Value = '  +6.04%  '
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("D51").Value = '''11.12'
$ws.Range("E51").Value = '  +20.99%  '
